$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '37.047.85'
$ws.Range('E2').Value = '  +0.10%  '

# Row 3
$ws.Range('D3').Value = '2.043.14'

# Row 4
$ws.Range('E4').Value = '  -0.16%  '

# Row 5
$ws.Range('D5').Value = '247.19'
$ws.Range('E5').Value = '  -1.36%  '

# Row 6
$ws.Range('E6').Value = '  -0.23%  '

# Row 8
$ws.Range('D8').Value = '56.31'
$ws.Range('E8').Value = '  +0.22%  '

# Row 9
$ws.Range('E9').Value = '  -0.33%  '

# Row 10
$ws.Range('D10').Value = '0.0779'
$ws.Range('E10').Value = '  -0.31%  '

# Row 11
$ws.Range('E11').Value = '  +0.44%  '

# Row 12
$ws.Range('D12').Value = '15.95'
$ws.Range('E12').Value = '  -2.95%  '

# Row 13
$ws.Range('E13').Value = '  +12.00%  '

# Row 14
$ws.Range('D14').Value = '2.339.86'
$ws.Range('E14').Value = '  -0.09%  '

# Row 15
$ws.Range('E15').Value = '  +2.46%  '

# Row 16
$ws.Range('D16').Value = '2.044.84'
$ws.Range('E16').Value = '  +0.10%  '

# Row 17
$ws.Range('D17').Value = '18.72'
$ws.Range('E17').Value = '  +12.15%  '

# Row 18
$ws.Range('D18').Value = '37.093.06'
$ws.Range('E18').Value = '  +0.40%  '

# Row 19
$ws.Range('D19').Value = '74.78'
$ws.Range('E19').Value = '  -0.14%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0890'
$ws.Range('E20').Value = '  -1.49%  '

# Row 21
$ws.Range('D21').Value = '5.40'
$ws.Range('E21').Value = '  +0.43%  '

# Row 22
$ws.Range('D22').Value = '236.67'
$ws.Range('E22').Value = '  +0.07%  '

# Row 23
$ws.Range('E23').Value = '  -0.13%  '

# Row 24
$ws.Range('D24').Value = '2.48'
$ws.Range('E24').Value = '  +4.77%  '

# Row 25
$ws.Range('D25').Value = '170.94'
$ws.Range('E25').Value = '  +1.15%  '

# Row 26
$ws.Range('E26').Value = '  +3.27%  '

# Row 27
$ws.Range('E27').Value = '  -8.58%  '

# Row 28
$ws.Range('D28').Value = '20.08'
$ws.Range('E28').Value = '  -0.13%  '

# Row 29
$ws.Range('E29').Value = '  -0.46%  '

# Row 30
$ws.Range('D30').Value = '5.11'
$ws.Range('E30').Value = '  +8.85%  '

# Row 31
$ws.Range('E31').Value = '  +0.66%  '

# Row 32
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '4.64'
$ws.Range('E32').Value = '  +4.71%  '

# Row 33
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.0621'
$ws.Range('E33').Value = '  +0.50%  '

# Row 34
$ws.Range('E34').Value = '  -0.11%  '

# Row 35
$ws.Range('E35').Value = '  -1.42%  '

# Row 36
$ws.Range('D36').Value = '1.86'
$ws.Range('E36').Value = '  +5.38%  '

# Row 37
$ws.Range('E37').Value = '  +1.41%  '

# Row 38
$ws.Range('E38').Value = '  -0.98%  '

# Row 39
$ws.Range('E39').Value = '  +9.28%  '

# Row 40
$ws.Range('E40').Value = '  +8.38%  '

# Row 41
$ws.Range('D41').Value = '0.0993'
$ws.Range('E41').Value = '  -8.98%  '

# Row 42
$ws.Range('E42').Value = '  +0.22%  '

# Row 43
$ws.Range('E43').Value = '  +1.56%  '

# Row 44
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '98.15'
$ws.Range('E44').Value = '  +1.69%  '

# Row 45
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '17.15'
$ws.Range('E45').Value = '  -2.41%  '

# Row 46
$ws.Range('D46').Value = '2.39'
$ws.Range('E46').Value = '  -3.62%  '

# Row 47
$ws.Range('D47').Value = '1.282.13'
$ws.Range('E47').Value = '  +0.17%  '

# Row 48
$ws.Range('E48').Value = '  -1.40%  '

# Row 49
$ws.Range('E49').Value = '  +0.98%  '

# Row 50
$ws.Range('D50').Value = '2.224.22'
$ws.Range('E50').Value = '  -0.18%  '

# Row 51
$ws.Range('D51').Value = '44.47'
$ws.Range('E51').Value = '  +1.97%  '
